$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 12933.333
$ws.Range("J26").Value = 12933.333
$ws.Range("L26").Value = 12933.333
$ws.Range("N26").Value = -13621.333
$ws.Range("H106").Value = 3720.2632
$ws.Range("I106").Value = 3605.3125
$ws.Range("J106").Value = 4333.3335
$ws.Range("K106").Value = 3605.3125
$ws.Range("L106").Value = 4333.3335
$ws.Range("M106").Value = -2974.3125
$ws.Range("N106").Value = -5595.3335
$ws.Range("H112").Value = 1989.2
$ws.Range("J112").Value = 2498.8572
$ws.Range("L112").Value = 7496.571599999999
$ws.Range("N112").Value = -9712.571599999999
$ws.Range("H123").Value = 33445.715
$ws.Range("I123").Value = 30000
$ws.Range("J123").Value = 33710.77
$ws.Range("K123").Value = 30000
$ws.Range("L123").Value = 33710.77
$ws.Range("M123").Value = -25100
$ws.Range("N123").Value = -43510.77
$ws.Range("H129").Value = 1085.9056
$ws.Range("J129").Value = 1133.08
$ws.Range("L129").Value = 3399.24
$ws.Range("N129").Value = -13399.24
$ws.Range("H131").Value = 3522.2415
$ws.Range("I131").Value = 466.15384
$ws.Range("J131").Value = 6005.3125
$ws.Range("K131").Value = 1398.46152
$ws.Range("L131").Value = 18015.9375
$ws.Range("M131").Value = 3641.53848
$ws.Range("N131").Value = -28095.9375
$ws.Range("H137").Value = 580869.9399999999
$ws.Range("I137").Value = 1926.3903
$ws.Range("J137").Value = 3971824.8
$ws.Range("K137").Value = 5779.1709
$ws.Range("L137").Value = 11915474.4
$ws.Range("M137").Value = -3229.1709
$ws.Range("N137").Value = -11920574.4
$ws.Range("H138").Value = 5010.036
$ws.Range("I138").Value = 2562.8572
$ws.Range("J138").Value = 5845.6587
$ws.Range("K138").Value = 7688.571599999999
$ws.Range("L138").Value = 17536.9761
$ws.Range("M138").Value = -2548.571599999999
$ws.Range("N138").Value = -27816.9761
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19517.592
$ws.Range("I32").Value = 20806.582
$ws.Range("J32").Value = 13072.637
$ws.Range("K32").Value = 20806.582
$ws.Range("L32").Value = 13072.637
$ws.Range("M32").Value = -20519.582
$ws.Range("N32").Value = -13646.637
$ws.Range("H45").Value = 1546.1212
$ws.Range("I45").Value = 1526.138
$ws.Range("K45").Value = 1526.138
$ws.Range("M45").Value = -1149.138
$ws.Range("H61").Value = 5719.2666
$ws.Range("I61").Value = 3816.5881
$ws.Range("J61").Value = 11600.272
$ws.Range("K61").Value = 3816.5881
$ws.Range("L61").Value = 11600.272
$ws.Range("M61").Value = -3604.5881
$ws.Range("N61").Value = -12024.272
$ws.Range("H102").Value = 3459
$ws.Range("I102").Value = 2898.3333
$ws.Range("K102").Value = 2898.3333
$ws.Range("M102").Value = -1276.3333
$ws.Range("H132").Value = 1645.4524
$ws.Range("I132").Value = 1016.65717
$ws.Range("J132").Value = 4789.4287
$ws.Range("K132").Value = 3049.97151
$ws.Range("L132").Value = 14368.2861
$ws.Range("M132").Value = -519.9715099999999
$ws.Range("N132").Value = -19428.2861
$ws.Range("H136").Value = 5719.2666
$ws.Range("I136").Value = 3816.5881
$ws.Range("J136").Value = 11600.272
$ws.Range("K136").Value = 11449.7643
$ws.Range("L136").Value = 34800.81600000001
$ws.Range("M136").Value = -8899.764299999999
$ws.Range("N136").Value = -39900.81600000001
$ws.Range("H138").Value = 38711
$ws.Range("J138").Value = 38711
$ws.Range("L138").Value = 38711
$ws.Range("N138").Value = -48991
$ws.Range("H139").Value = 36299.11
$ws.Range("J139").Value = 36299.11
$ws.Range("L139").Value = 36299.11
$ws.Range("N139").Value = -46579.11
$ws.Range("H140").Value = 41117.668
$ws.Range("J140").Value = 41117.668
$ws.Range("L140").Value = 41117.668
$ws.Range("N140").Value = -51477.668
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 984.8461
$ws.Range("I20").Value = 1011.2727
$ws.Range("J20").Value = 839.5
$ws.Range("K20").Value = 1011.2727
$ws.Range("L20").Value = 839.5
$ws.Range("M20").Value = -764.2727
$ws.Range("N20").Value = -1333.5
$ws.Range("H80").Value = 171.22223
$ws.Range("J80").Value = 175.41176
$ws.Range("L80").Value = 175.41176
$ws.Range("N80").Value = -2171.41176
$ws.Range("H83").Value = 171.22223
$ws.Range("J83").Value = 175.41176
$ws.Range("L83").Value = 877.0587999999999
$ws.Range("N83").Value = -10861.0588
$ws.Range("H99").Value = 983.3333
$ws.Range("I99").Value = 980
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 980
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 518
$ws.Range("N99").Value = -3996
$ws.Range("H105").Value = 7269.231
$ws.Range("I105").Value = 10000
$ws.Range("J105").Value = 6055.5557
$ws.Range("K105").Value = 10000
$ws.Range("L105").Value = 6055.5557
$ws.Range("M105").Value = -8253
$ws.Range("N105").Value = -9549.555700000001
$ws.Range("H131").Value = 24345.455
$ws.Range("J131").Value = 24345.455
$ws.Range("L131").Value = 24345.455
$ws.Range("N131").Value = -34425.455
$ws.Range("H134").Value = 1903.8108
$ws.Range("I134").Value = 1897.2858
$ws.Range("J134").Value = 1924.1111
$ws.Range("K134").Value = 5691.857400000001
$ws.Range("L134").Value = 5772.3333
$ws.Range("M134").Value = -3156.857400000001
$ws.Range("N134").Value = -10842.3333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 14942.75
$ws.Range("J50").Value = 14942.75
$ws.Range("L50").Value = 14942.75
$ws.Range("N50").Value = -16192.75
$ws.Range("H59").Value = 14999
$ws.Range("J59").Value = 14999
$ws.Range("L59").Value = 14999
$ws.Range("N59").Value = -17289
$ws.Range("H134").Value = 3371.2
$ws.Range("I134").Value = 1978.5883
$ws.Range("J134").Value = 4400.522
$ws.Range("K134").Value = 5935.7649
$ws.Range("L134").Value = 13201.566
$ws.Range("M134").Value = -3400.7649
$ws.Range("N134").Value = -18271.566
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 431
$ws.Range("J6").Value = 668
$ws.Range("L6").Value = 2004
$ws.Range("N6").Value = -2230
$ws.Range("H51").Value = 9345.166999999999
$ws.Range("I51").Value = 706.6667
$ws.Range("K51").Value = 2120.0001
$ws.Range("M51").Value = -1660.0001
$ws.Range("H63").Value = 3148
$ws.Range("I63").Value = 1500.2
$ws.Range("J63").Value = 3736.5
$ws.Range("K63").Value = 4500.6
$ws.Range("L63").Value = 11209.5
$ws.Range("M63").Value = -3751.6
$ws.Range("N63").Value = -12707.5
$ws.Range("H66").Value = 3148
$ws.Range("I66").Value = 1500.2
$ws.Range("J66").Value = 3736.5
$ws.Range("K66").Value = 13501.8
$ws.Range("L66").Value = 33628.5
$ws.Range("M66").Value = -9757.800000000001
$ws.Range("N66").Value = -41116.5
$ws.Range("H70").Value = 2746.5454
$ws.Range("I70").Value = 1242.4
$ws.Range("K70").Value = 3727.2
$ws.Range("M70").Value = -3412.2
$ws.Range("H73").Value = 2746.5454
$ws.Range("I73").Value = 1242.4
$ws.Range("K73").Value = 3727.2
$ws.Range("M73").Value = -2635.2
$ws.Range("H86").Value = 1266.25
$ws.Range("I86").Value = 1343.5555
$ws.Range("J86").Value = 1034.3334
$ws.Range("K86").Value = 4030.6665
$ws.Range("L86").Value = 3103.0002
$ws.Range("M86").Value = -2844.6665
$ws.Range("N86").Value = -5475.0002
$ws.Range("H89").Value = 1266.25
$ws.Range("I89").Value = 1343.5555
$ws.Range("J89").Value = 1034.3334
$ws.Range("K89").Value = 12091.9995
$ws.Range("L89").Value = 9309.000599999999
$ws.Range("M89").Value = -6163.9995
$ws.Range("N89").Value = -21165.0006
$ws.Range("H131").Value = 34444.45
$ws.Range("I131").Value = 1450.9412
$ws.Range("J131").Value = 81185.25
$ws.Range("K131").Value = 4352.8236
$ws.Range("L131").Value = 243555.75
$ws.Range("M131").Value = 687.1764000000003
$ws.Range("N131").Value = -253635.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 24920
$ws.Range("H97").Value = 14752.308
$ws.Range("I97").Value = 26140
$ws.Range("J97").Value = 1466.6666
$ws.Range("K97").Value = 26140
$ws.Range("L97").Value = 1466.6666
$ws.Range("M97").Value = -25644
$ws.Range("N97").Value = -2458.6666
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 944.25
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 790.8
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 790.8
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -1166.8
$ws.Range("H132").Value = 5175.6665
$ws.Range("I132").Value = 6469.6
$ws.Range("J132").Value = 3272.8235
$ws.Range("K132").Value = 19408.8
$ws.Range("L132").Value = 9818.470499999999
$ws.Range("M132").Value = -16878.8
$ws.Range("N132").Value = -14878.4705
$ws.Range("H136").Value = 4710.4683
$ws.Range("I136").Value = 2627.7693
$ws.Range("J136").Value = 7289.048
$ws.Range("K136").Value = 7883.3079
$ws.Range("L136").Value = 21867.144
$ws.Range("M136").Value = -5333.3079
$ws.Range("N136").Value = -26967.144

Write-Host "Updated $($wb.ActiveSheet) cells"